$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows (row 1 = header, rows 2-12 = data)
# Columns: A = Variable, B = Lambda, C = P, D = P (hochberg-corrected)
$data = @(
    @("BIO1",      "0.86001899999999998",  "1.13812E-48",            "1.026E-47"),
    @("BIO2",      "0.67059999999999997",  "4.2750000000000001E-122","4.28E-121"),
    @("BIO3",      "0.88739999999999997",  "1.8553E-134",            "2.0459999999999999E-133"),
    @("BIO4",      "0.75571699999999997",  "2.7908199999999999E-42", "2.2320000000000001E-41"),
    @("BIO7",      "0.68789599999999995",  "5.7569500000000002E-30", "2.8800000000000002E-29"),
    @("BIO12",     "0.585947",             "2.00525E-13",            "2.01E-13"),
    @("BIO15",     "0.71018000000000003",  "1.6742200000000001E-26", "5.0100000000000002E-26"),
    @("elevation", "0.85695200000000005",  "1.4637700000000001E-35", "1.022E-34"),
    @("nitrogen",  "0.54563399999999995",  "1.9351200000000001E-29", "7.7599999999999995E-29"),
    @("pH",        "0.31294300000000003",  "3.5869299999999998E-24", "7.1799999999999995E-24"),
    @("aridity",   "0.63764799999999999",  "8.4531000000000004E-31", "5.0700000000000003E-30")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = [double]$r[1]
    $ws.Cells.Item($row, 3).Value = [double]$r[2]
    $ws.Cells.Item($row, 3).NumberFormat = "0.00E+00"
    $ws.Cells.Item($row, 4).Value = [double]$r[3]
    $ws.Cells.Item($row, 4).NumberFormat = "0.00E+00"
    $row = $row + 1
}

# Header row (added last so shared-string order matches source order)
$ws.Range("D1").Value = "P (hochberg-corrected)"

# Clear any leftover rows below row 12 (old data went up to row 20)
$ws.Range("A13:D20").Clear()

# Selection / view update per diff
$ws.Range("D2:D13").Select()
